# Auto-generated edit script applying odds/time updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value2 = 1.07
$ws.Range("O2").Value2 = 1.36
$ws.Range("M3").Value2 = 1.1
$ws.Range("O3").Value2 = 1.4
$ws.Range("U3").Value2 = 1.95
$ws.Range("V3").Value2 = 1.8
$ws.Range("N4").Value2 = 8
$ws.Range("AC4").Value2 = 8
$ws.Range("AD4").Value2 = 7
$ws.Range("AE4").Value2 = 21
$ws.Range("AI4").Value2 = 21
$ws.Range("AJ4").Value2 = 15
$ws.Range("G5").Value2 = 2.5
$ws.Range("I5").Value2 = 3.2
$ws.Range("M5").Value2 = 1.11
$ws.Range("N5").Value2 = 6.5
$ws.Range("Z5").Value2 = 23
$ws.Range("AA5").Value2 = 23
$ws.Range("AI5").Value2 = 15
$ws.Range("AJ5").Value2 = 13
$ws.Range("G7").Value2 = 2.5
$ws.Range("I7").Value2 = 2.7
$ws.Range("J7").Value2 = 3.2
$ws.Range("L7").Value2 = 3.4
$ws.Range("Q7").Value2 = 1.98
$ws.Range("R7").Value2 = 1.88
$ws.Range("W7").Value2 = 8.5
$ws.Range("X7").Value2 = 12
$ws.Range("Z7").Value2 = 23
$ws.Range("AK7").Value2 = 29
$ws.Range("AN7").Value2 = 4.5
$ws.Range("AQ7").Value2 = 41
$ws.Range("G8").Value2 = 4.33
$ws.Range("H8").Value2 = 3.5
$ws.Range("I8").Value2 = 1.7
$ws.Range("K8").Value2 = 2.25
$ws.Range("L8").Value2 = 2.38
$ws.Range("Q8").Value2 = 1.9
$ws.Range("R8").Value2 = 1.95
$ws.Range("S8").Value2 = 1.36
$ws.Range("T8").Value2 = 3
$ws.Range("U8").Value2 = 1.8
$ws.Range("V8").Value2 = 1.91
$ws.Range("AA8").Value2 = 41
$ws.Range("AC8").Value2 = 11
$ws.Range("AK8").Value2 = 13
$ws.Range("AL8").Value2 = 13
$ws.Range("AN8").Value2 = 6.5
$ws.Range("AT8").Value2 = 3
$ws.Range("AX8").Value2 = 9
$ws.Range("AY8").Value2 = 19
$ws.Range("BB8").Value2 = 126
$ws.Range("Q9").Value2 = 1.75
$ws.Range("U9").Value2 = 1.67
$ws.Range("Q11").Value2 = 1.57
$ws.Range("R11").Value2 = 2.35
$ws.Range("Q12").Value2 = 1.67
$ws.Range("G15").Value2 = 1.67
$ws.Range("H15").Value2 = 3.9
$ws.Range("I15").Value2 = 4.5
$ws.Range("J15").Value2 = 2.25
$ws.Range("L15").Value2 = 4.75
$ws.Range("O15").Value2 = 1.18
$ws.Range("P15").Value2 = 4.5
$ws.Range("Q15").Value2 = 1.62
$ws.Range("R15").Value2 = 2.25
$ws.Range("U15").Value2 = 1.67
$ws.Range("V15").Value2 = 2.1
$ws.Range("X15").Value2 = 9
$ws.Range("Z15").Value2 = 13
$ws.Range("AI15").Value2 = 26
$ws.Range("AJ15").Value2 = 15
$ws.Range("AK15").Value2 = 51
$ws.Range("AL15").Value2 = 34
$ws.Range("AO15").Value2 = 8.5
$ws.Range("AQ15").Value2 = 26
$ws.Range("AW15").Value2 = 6.5
$ws.Range("AX15").Value2 = 23
$ws.Range("G16").Value2 = 1.5
$ws.Range("I16").Value2 = 5.75
$ws.Range("J16").Value2 = 2.05
$ws.Range("K16").Value2 = 2.38
$ws.Range("Q16").Value2 = 1.75
$ws.Range("R16").Value2 = 2.05
$ws.Range("S16").Value2 = 1.33
$ws.Range("T16").Value2 = 3.25
$ws.Range("U16").Value2 = 1.83
$ws.Range("V16").Value2 = 1.83
$ws.Range("X16").Value2 = 7.5
$ws.Range("AA16").Value2 = 12
$ws.Range("AD16").Value2 = 8.5
$ws.Range("AT16").Value2 = 3.25
$ws.Range("M17").Value2 = 1.01
$ws.Range("O17").Value2 = 1.08
$ws.Range("M18").Value2 = 1.03
$ws.Range("O18").Value2 = 1.25
$ws.Range("U18").Value2 = 1.73
$ws.Range("G19").Value2 = 1.7
$ws.Range("I19").Value2 = 4.33
$ws.Range("L19").Value2 = 4.75
$ws.Range("M19").Value2 = 1.02
$ws.Range("N19").Value2 = 15
$ws.Range("O19").Value2 = 1.17
$ws.Range("Q19").Value2 = 1.67
$ws.Range("U19").Value2 = 1.67
$ws.Range("Z19").Value2 = 13
$ws.Range("AK19").Value2 = 51
$ws.Range("AN19").Value2 = 3.75
$ws.Range("AO19").Value2 = 8.5
$ws.Range("AQ19").Value2 = 26
$ws.Range("AX19").Value2 = 23
$ws.Range("Q20").Value2 = 1.48
$ws.Range("U20").Value2 = 1.67
$ws.Range("Q21").Value2 = 1.33
$ws.Range("U21").Value2 = 1.53
$ws.Range("V21").Value2 = 2.38
$ws.Range("C22").Value = "15:15"
$ws.Range("Q22").Value2 = 1.73
$ws.Range("R22").Value2 = 2.08
$ws.Range("U22").Value2 = 1.57
$ws.Range("U24").Value2 = 1.87
$ws.Range("V24").Value2 = 1.77
$ws.Range("U25").Value2 = 1.77
$ws.Range("V25").Value2 = 1.92
$ws.Range("U26").Value2 = 1.58
$ws.Range("U28").Value2 = 1.77
$ws.Range("V28").Value2 = 1.87
$ws.Range("M30").Value2 = 1.05
$ws.Range("O30").Value2 = 1.27
$ws.Range("U30").Value2 = 1.77
$ws.Range("V30").Value2 = 1.87
$ws.Range("G32").Value2 = 2.05
$ws.Range("I32").Value2 = 3.25
$ws.Range("J32").Value2 = 2.63
$ws.Range("L32").Value2 = 3.6
$ws.Range("M32").Value2 = 1.02
$ws.Range("O32").Value2 = 1.15
$ws.Range("P32").Value2 = 4.5
$ws.Range("Q32").Value2 = 1.62
$ws.Range("R32").Value2 = 2.25
$ws.Range("U32").Value2 = 1.54
$ws.Range("V32").Value2 = 2.25
$ws.Range("W32").Value2 = 10
$ws.Range("X32").Value2 = 12
$ws.Range("Y32").Value2 = 9
$ws.Range("Z32").Value2 = 19
$ws.Range("AA32").Value2 = 15
$ws.Range("AD32").Value2 = 7.5
$ws.Range("AE32").Value2 = 12
$ws.Range("AF32").Value2 = 41
$ws.Range("AG32").Value2 = 126
$ws.Range("AI32").Value2 = 19
$ws.Range("AJ32").Value2 = 12
$ws.Range("AK32").Value2 = 34
$ws.Range("AL32").Value2 = 23
$ws.Range("AM32").Value2 = 26
$ws.Range("AN32").Value2 = 4.33
$ws.Range("AO32").Value2 = 11
$ws.Range("AP32").Value2 = 17
$ws.Range("AQ32").Value2 = 34
$ws.Range("AW32").Value2 = 5.5
$ws.Range("AX32").Value2 = 17
$ws.Range("AZ32").Value2 = 51
$ws.Range("BA32").Value2 = 67
$ws.Range("BB32").Value2 = 126
$ws.Range("BC32").Value2 = 401
$ws.Range("G33").Value2 = 2.9
$ws.Range("I33").Value2 = 2.55
$ws.Range("J33").Value2 = 3.6
$ws.Range("L33").Value2 = 3.25
$ws.Range("N33").Value2 = 7.5
$ws.Range("R33").Value2 = 1.62
$ws.Range("U33").Value2 = 1.92
$ws.Range("V33").Value2 = 1.77
$ws.Range("W33").Value2 = 8
$ws.Range("X33").Value2 = 13
$ws.Range("Y33").Value2 = 11
$ws.Range("Z33").Value2 = 29
$ws.Range("AA33").Value2 = 26
$ws.Range("AH33").Value2 = 7.5
$ws.Range("AL33").Value2 = 23
$ws.Range("AN33").Value2 = 4.75
$ws.Range("AO33").Value2 = 17
$ws.Range("AW33").Value2 = 4.5
$ws.Range("AX33").Value2 = 15
$ws.Range("BA33").Value2 = 81
$ws.Range("G34").Value2 = 1.6
$ws.Range("N34").Value2 = 13
$ws.Range("Q34").Value2 = 1.75
$ws.Range("R34").Value2 = 2.05
$ws.Range("S34").Value2 = 1.33
$ws.Range("T34").Value2 = 3.25
$ws.Range("U34").Value2 = 1.77
$ws.Range("V34").Value2 = 1.92
$ws.Range("W34").Value2 = 7.5
$ws.Range("X34").Value2 = 8
$ws.Range("AB34").Value2 = 23
$ws.Range("AC34").Value2 = 13
$ws.Range("AG34").Value2 = 251
$ws.Range("AT34").Value2 = 3.25
$ws.Range("AU34").Value2 = 8
$ws.Range("AX34").Value2 = 26
$ws.Range("G35").Value2 = 2.38
$ws.Range("I35").Value2 = 3.2
$ws.Range("J35").Value2 = 3.1
$ws.Range("L35").Value2 = 3.75
$ws.Range("R35").Value2 = 1.65
$ws.Range("U35").Value2 = 1.8
$ws.Range("V35").Value2 = 1.8
$ws.Range("W35").Value2 = 7
$ws.Range("X35").Value2 = 11
$ws.Range("Y35").Value2 = 9.5
$ws.Range("Z35").Value2 = 21
$ws.Range("AH35").Value2 = 9
$ws.Range("AI35").Value2 = 15
$ws.Range("AJ35").Value2 = 12
$ws.Range("AK35").Value2 = 34
$ws.Range("AL35").Value2 = 29
$ws.Range("AM35").Value2 = 41
$ws.Range("AO35").Value2 = 13
$ws.Range("AQ35").Value2 = 41
$ws.Range("AW35").Value2 = 5
$ws.Range("G36").Value2 = 2.15
$ws.Range("I36").Value2 = 3.1
$ws.Range("M36").Value2 = 1.06
$ws.Range("N36").Value2 = 10
$ws.Range("Y36").Value2 = 9
$ws.Range("AA36").Value2 = 17
$ws.Range("AI36").Value2 = 17
$ws.Range("AJ36").Value2 = 12
